# Add 2022-Q3 data:
#  1. Insert a new worksheet named "2022-Q3" right after the "总计" sheet,
#     populated with the quarterly per-fund holdings table.
#  2. Update the "总计" (summary) sheet: add a new top data row for
#     2022-Q3 and append the row that got pushed off the bottom
#     (2020-Q4), so the whole table now spans 8 data rows (A1:D9).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"

# Header row. Match the style used for header cells on the other
# quarter sheets (bold, centered, bordered) by copying the format from
# the "总计" sheet's own header cell, then fill in the header text.
$totalSheet.Range("B1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3Sheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# Fund holdings rows (code, name, scale, stock position, position ratio,
# market value (CNY bn), position rank). The numeric-looking columns
# (scale/positions/ratio/value) are stored as text, matching the source
# sheets; only the index (A) and rank (H) columns are real numbers.
$rows = @(
    @("011463", "长城量化精选股票C", "1.53", "94.12", "6.73", "0.1030", 8),
    @("006926", "长城量化精选股票A", "1.06", "94.12", "6.73", "0.0713", 8),
    @("200016", "长城稳健成长灵活配置混合", "0.73", "79.86", "2.74", "0.0200", 6),
    @("002159", "东吴国企改革主题灵活配置混合A", "0.17", "92.20", "6.02", "0.0102", 10),
    @("015741", "东财品质生活优选混合A", "0.10", "93.68", "8.42", "0.0084", 4),
    @("012615", "东吴国企改革主题灵活配置混合C", "0.13", "92.20", "6.02", "0.0078", 10),
    @("015742", "东财品质生活优选混合C", "0.02", "93.68", "8.42", "0.0017", 4)
)

$q3Sheet.Range("B2:G8").NumberFormat = "@"

# Column A (row index) uses the same bordered/centered style as the
# header and as column A on the other quarter sheets.
$totalSheet.Range("A2").Copy()
$q3Sheet.Range("A2:A8").PasteSpecial(-4122)  # xlPasteFormats

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $q3Sheet.Cells.Item($r, 1).Value = $i
    $q3Sheet.Cells.Item($r, 2).Value = $rows[$i][0]
    $q3Sheet.Cells.Item($r, 3).Value = $rows[$i][1]
    $q3Sheet.Cells.Item($r, 4).Value = $rows[$i][2]
    $q3Sheet.Cells.Item($r, 5).Value = $rows[$i][3]
    $q3Sheet.Cells.Item($r, 6).Value = $rows[$i][4]
    $q3Sheet.Cells.Item($r, 7).Value = $rows[$i][5]
    $q3Sheet.Cells.Item($r, 8).Value = $rows[$i][6]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet with the new 2022-Q3 row, shifting the
#    existing quarters down and re-adding 2020-Q4 at the new last row.
# ---------------------------------------------------------------------
$totalSheet.Range("A8").Copy()
$totalSheet.Range("A9").PasteSpecial(-4122)  # xlPasteFormats

$summary = @(
    @("2022-Q3", 7, 0.22),
    @("2022-Q2", 23, 19.72),
    @("2022-Q1", 13, 19.52),
    @("2021-Q4", 28, 28.23),
    @("2021-Q3", 66, 47.73),
    @("2021-Q2", 40, 20.86),
    @("2021-Q1", 26, 4.79),
    @("2020-Q4", 24, 16.4)
)

for ($i = 0; $i -lt $summary.Length; $i++) {
    $r = $i + 2
    $totalSheet.Cells.Item($r, 1).Value = $i
    $totalSheet.Cells.Item($r, 2).Value = $summary[$i][0]
    $totalSheet.Cells.Item($r, 3).Value = $summary[$i][1]
    $totalSheet.Cells.Item($r, 4).Value = $summary[$i][2]
}

# Keep "总计" as the active sheet, same as before the edit.
$totalSheet.Activate()
